# Updates the cryptocurrency price/volume table to reflect the latest
# scrape, per "Updated cryptos list on Mon Feb 12 16:29:00 UTC 2024 with
# GitHub Actions".
#
# Values are written as text (matching the source data, which stores
# prices/percentages as text, e.g. "49.691.21" or "  +3.30%  "). Forcing
# the cell's number format to Text ("@") before the assignment prevents
# Excel from re-interpreting numeric-looking strings (like "2.00" or
# "321.71") as numbers, which would otherwise drop trailing zeros or
# introduce floating point noise. Resetting the style back to "Normal"
# afterwards avoids leaving a stray text-format style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '49.691.21' },
    @{ Cell = 'E2'; Value = '  +3.30%  ' },
    @{ Cell = 'D3'; Value = '2.548.25' },
    @{ Cell = 'E3'; Value = '  +1.66%  ' },
    @{ Cell = 'D4'; Value = '0.998' },
    @{ Cell = 'E4'; Value = '  -0.08%  ' },
    @{ Cell = 'D5'; Value = '321.71' },
    @{ Cell = 'E5'; Value = '  +0.14%  ' },
    @{ Cell = 'D6'; Value = '108.51' },
    @{ Cell = 'E6'; Value = '  +0.11%  ' },
    @{ Cell = 'D7'; Value = '0.526' },
    @{ Cell = 'E7'; Value = '  -0.20%  ' },
    @{ Cell = 'D8'; Value = '0.998' },
    @{ Cell = 'E8'; Value = '  -0.07%  ' },
    @{ Cell = 'E9'; Value = '  +2.66%  ' },
    @{ Cell = 'D10'; Value = '40.14' },
    @{ Cell = 'E10'; Value = '  +0.83%  ' },
    @{ Cell = 'D11'; Value = '20.28' },
    @{ Cell = 'E11'; Value = '  +0.61%  ' },
    @{ Cell = 'D12'; Value = '0.0813' },
    @{ Cell = 'E12'; Value = '  -0.47%  ' },
    @{ Cell = 'E13'; Value = '  +0.38%  ' },
    @{ Cell = 'E14'; Value = '  +0.64%  ' },
    @{ Cell = 'D15'; Value = '2.945.90' },
    @{ Cell = 'E15'; Value = '  +1.63%  ' },
    @{ Cell = 'D16'; Value = '2.586.91' },
    @{ Cell = 'E16'; Value = '  +3.36%  ' },
    @{ Cell = 'D17'; Value = '0.856' },
    @{ Cell = 'E17'; Value = '  +1.33%  ' },
    @{ Cell = 'D18'; Value = '49.496.28' },
    @{ Cell = 'E18'; Value = '  +3.22%  ' },
    @{ Cell = 'B19'; Value = 'ImmutableX' },
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Cell = 'D19'; Value = '3.03' },
    @{ Cell = 'E19'; Value = '  +11.05%  ' },
    @{ Cell = 'B20'; Value = 'InternetComputer(DFINITY)' },
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' },
    @{ Cell = 'D20'; Value = '13.20' },
    @{ Cell = 'E20'; Value = '  +0.38%  ' },
    @{ Cell = 'D21'; Value = '6.67' },
    @{ Cell = 'E21'; Value = '  +0.45%  ' },
    @{ Cell = 'E22'; Value = '  -0.05%  ' },
    @{ Cell = 'D23'; Value = '286.96' },
    @{ Cell = 'E23'; Value = '  +3.90%  ' },
    @{ Cell = 'D24'; Value = '71.90' },
    @{ Cell = 'E24'; Value = '  -0.24%  ' },
    @{ Cell = 'D25'; Value = '2.52' },
    @{ Cell = 'E25'; Value = '  -1.37%  ' },
    @{ Cell = 'D26'; Value = '26.36' },
    @{ Cell = 'E26'; Value = '  +1.83%  ' },
    @{ Cell = 'E27'; Value = '  -0.17%  ' },
    @{ Cell = 'E28'; Value = '  -1.24%  ' },
    @{ Cell = 'E29'; Value = '  +3.96%  ' },
    @{ Cell = 'D30'; Value = '9.82' },
    @{ Cell = 'E30'; Value = '  -2.14%  ' },
    @{ Cell = 'D31'; Value = '35.21' },
    @{ Cell = 'E31'; Value = '  -0.66%  ' },
    @{ Cell = 'D32'; Value = '49.44' },
    @{ Cell = 'E32'; Value = '  -0.03%  ' },
    @{ Cell = 'D34'; Value = '5.35' },
    @{ Cell = 'E34'; Value = '  +0.16%  ' },
    @{ Cell = 'E35'; Value = '  -0.16%  ' },
    @{ Cell = 'D36'; Value = '0.0784' },
    @{ Cell = 'E36'; Value = '  -0.09%  ' },
    @{ Cell = 'D37'; Value = '2.00' },
    @{ Cell = 'E37'; Value = '  +2.45%  ' },
    @{ Cell = 'E38'; Value = '  +1.28%  ' },
    @{ Cell = 'D39'; Value = '2.99' },
    @{ Cell = 'E39'; Value = '  +0.93%  ' },
    @{ Cell = 'E40'; Value = '  +0.14%  ' },
    @{ Cell = 'B41'; Value = 'EnergySwap' },
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'D41'; Value = '22.35' },
    @{ Cell = 'E41'; Value = '  +2.79%  ' },
    @{ Cell = 'B42'; Value = 'Monero' },
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = 'D42'; Value = '120.64' },
    @{ Cell = 'E42'; Value = '  -1.85%  ' },
    @{ Cell = 'B43'; Value = 'WEMIXToken' },
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' },
    @{ Cell = 'D43'; Value = '2.22' },
    @{ Cell = 'E43'; Value = '  +0.14%  ' },
    @{ Cell = 'D44'; Value = '0.0309' },
    @{ Cell = 'E44'; Value = '  +1.22%  ' },
    @{ Cell = 'E45'; Value = '  +4.59%  ' },
    @{ Cell = 'D46'; Value = '2.011.04' },
    @{ Cell = 'E46'; Value = '  +0.48%  ' },
    @{ Cell = 'D47'; Value = '2.01' },
    @{ Cell = 'E47'; Value = '  +8.46%  ' },
    @{ Cell = 'E48'; Value = '  +7.72%  ' },
    @{ Cell = 'E49'; Value = '  +0.18%  ' },
    @{ Cell = 'E50'; Value = '  +2.46%  ' },
    @{ Cell = 'D51'; Value = '81.41' },
    @{ Cell = 'E51'; Value = '  +2.03%  ' }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.Style = "Normal"
}
